# "Add files via upload" - the uploaded CV asset for color #3 ("Magenta-Indigo")
# was renamed/re-uploaded on GitHub, so its Icon_url cell needs to point at the
# new blob path (new commit hash + underscore-style filename instead of the
# URL-encoded space).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colors_table")

$ws.Range("E4").Value = "https://github.com/Ing-Aladar-Dukay/CV_Dukay/blob/622c97a85b9033e5c2bd2c66931c4a8961f7bb1d/03%20Colors%20icons/color_03.png"
